$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Save" header in H1, matching the formatting (bold header style)
# used by the other header cells (copy format from the neighboring "sum"
# header in G1 so the same shared cellXfs style gets reused).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the Save column values for rows 2-5 (plain numeric, unstyled - same
# as the other data cells in the row)
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
